$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# Row 129 (ALC)
$ws_ALC.Range("H129").Value = 825.7
$ws_ALC.Range("J129").Value = 865.44086
$ws_ALC.Range("L129").Value = 2596.32258
$ws_ALC.Range("N129").Value = -12596.32258

# Row 139 (ALC)
$ws_ALC.Range("H139").Value = 38360
$ws_ALC.Range("J139").Value = 38360
$ws_ALC.Range("L139").Value = 38360
$ws_ALC.Range("N139").Value = -48640

$ws_ARM = $wb.Worksheets.Item("ARM")
# Row 88 (ARM)
$ws_ARM.Range("H88").Value = 2951.1667
$ws_ARM.Range("I88").Value = 2000
$ws_ARM.Range("J88").Value = 3426.75
$ws_ARM.Range("K88").Value = 2000
$ws_ARM.Range("L88").Value = 3426.75
$ws_ARM.Range("M88").Value = -1594
$ws_ARM.Range("N88").Value = -4238.75

# Row 91 (ARM)
$ws_ARM.Range("H91").Value = 2951.1667
$ws_ARM.Range("I91").Value = 2000
$ws_ARM.Range("J91").Value = 3426.75
$ws_ARM.Range("K91").Value = 2000
$ws_ARM.Range("L91").Value = 3426.75
$ws_ARM.Range("M91").Value = -596
$ws_ARM.Range("N91").Value = -6234.75

$ws_BSM = $wb.Worksheets.Item("BSM")
# Row 86 (BSM)
$ws_BSM.Range("H86").Value = 1166.1875
$ws_BSM.Range("I86").Value = 1013.8182
$ws_BSM.Range("J86").Value = 1501.4
$ws_BSM.Range("K86").Value = 1013.8182
$ws_BSM.Range("L86").Value = 1501.4
$ws_BSM.Range("M86").Value = 109.1818
$ws_BSM.Range("N86").Value = -3747.4

# Row 89 (BSM)
$ws_BSM.Range("H89").Value = 1166.1875
$ws_BSM.Range("I89").Value = 1013.8182
$ws_BSM.Range("J89").Value = 1501.4
$ws_BSM.Range("K89").Value = 5069.091
$ws_BSM.Range("L89").Value = 7507
$ws_BSM.Range("M89").Value = 546.9089999999997
$ws_BSM.Range("N89").Value = -18739

# Row 95 (BSM)
$ws_BSM.Range("H95").Value = 33555.555
$ws_BSM.Range("J95").Value = 33555.555
$ws_BSM.Range("L95").Value = 33555.555
$ws_BSM.Range("N95").Value = -39047.555

# Row 105 (BSM)
$ws_BSM.Range("H105").Value = 1742.1807
$ws_BSM.Range("I105").Value = 1744.8572
$ws_BSM.Range("J105").Value = 1707.8334
$ws_BSM.Range("K105").Value = 1744.8572
$ws_BSM.Range("L105").Value = 1707.8334
$ws_BSM.Range("M105").Value = 2.142800000000079
$ws_BSM.Range("N105").Value = -5201.8334

# Row 134 (BSM)
$ws_BSM.Range("H134").Value = 2963.0322
$ws_BSM.Range("I134").Value = 1628.7273
$ws_BSM.Range("K134").Value = 4886.1819
$ws_BSM.Range("M134").Value = -2351.1819

$ws_CRP = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws_CRP.Range("H31").Value = 5537.773
$ws_CRP.Range("I31").Value = 2480.2068
$ws_CRP.Range("K31").Value = 2480.2068
$ws_CRP.Range("M31").Value = -2185.2068

# Row 34 (CRP)
$ws_CRP.Range("H34").Value = 5537.773
$ws_CRP.Range("I34").Value = 2480.2068
$ws_CRP.Range("K34").Value = 2480.2068
$ws_CRP.Range("M34").Value = -2278.2068

# Row 58 (CRP)
$ws_CRP.Range("H58").Value = 1915.3611
$ws_CRP.Range("I58").Value = 1449.8276
$ws_CRP.Range("J58").Value = 3844
$ws_CRP.Range("K58").Value = 1449.8276
$ws_CRP.Range("L58").Value = 3844
$ws_CRP.Range("M58").Value = -1246.8276
$ws_CRP.Range("N58").Value = -4250

# Row 136 (CRP)
$ws_CRP.Range("H136").Value = 1915.3611
$ws_CRP.Range("I136").Value = 1449.8276
$ws_CRP.Range("J136").Value = 3844
$ws_CRP.Range("K136").Value = 4349.4828
$ws_CRP.Range("L136").Value = 11532
$ws_CRP.Range("M136").Value = -1799.4828
$ws_CRP.Range("N136").Value = -16632

# Row 139 (CRP)
$ws_CRP.Range("H139").Value = 85011
$ws_CRP.Range("J139").Value = 85011
$ws_CRP.Range("L139").Value = 85011
$ws_CRP.Range("N139").Value = -95291

$ws_CUL = $wb.Worksheets.Item("CUL")
# Row 4 (CUL)
$ws_CUL.Range("H4").Value = 22142.857
$ws_CUL.Range("I4").Value = 100266.664
$ws_CUL.Range("J4").Value = 836.36365
$ws_CUL.Range("K4").Value = 300799.992
$ws_CUL.Range("L4").Value = 2509.09095
$ws_CUL.Range("M4").Value = -300687.992
$ws_CUL.Range("N4").Value = -2733.09095

# Row 34 (CUL)
$ws_CUL.Range("H34").Value = 10345.267
$ws_CUL.Range("J34").Value = 7619
$ws_CUL.Range("L34").Value = 22857
$ws_CUL.Range("N34").Value = -23025

# Row 60 (CUL)
$ws_CUL.Range("H60").Value = 16451
$ws_CUL.Range("I60").Value = 189.5
$ws_CUL.Range("K60").Value = 568.5
$ws_CUL.Range("M60").Value = -317.5

# Row 61 (CUL)
$ws_CUL.Range("H61").Value = 151.6923
$ws_CUL.Range("I61").Value = 122.2
$ws_CUL.Range("J61").Value = 250
$ws_CUL.Range("K61").Value = 366.6
$ws_CUL.Range("L61").Value = 750
$ws_CUL.Range("M61").Value = -151.6
$ws_CUL.Range("N61").Value = -1180

# Row 133 (CUL)
$ws_CUL.Range("H133").Value = 5570
$ws_CUL.Range("I133").Value = 5666.25
$ws_CUL.Range("K133").Value = 16998.75
$ws_CUL.Range("M133").Value = -11938.75

$ws_GSM = $wb.Worksheets.Item("GSM")
# Row 113 (GSM)
$ws_GSM.Range("H113").Value = 2355.9167
$ws_GSM.Range("J113").Value = 1972
$ws_GSM.Range("L113").Value = 1972
$ws_GSM.Range("N113").Value = -6312

# Row 139 (GSM)
$ws_GSM.Range("H139").Value = 0
$ws_GSM.Range("J139").Value = 0
$ws_GSM.Range("L139").Value = 0
$ws_GSM.Range("N139").ClearContents()

$ws_LTW = $wb.Worksheets.Item("LTW")
# Row 62 (LTW)
$ws_LTW.Range("H62").Value = 20142.334
$ws_LTW.Range("I62").Value = 12227
$ws_LTW.Range("J62").Value = 24100
$ws_LTW.Range("K62").Value = 12227
$ws_LTW.Range("L62").Value = 24100
$ws_LTW.Range("M62").Value = -11603
$ws_LTW.Range("N62").Value = -25348

# Row 65 (LTW)
$ws_LTW.Range("H65").Value = 20142.334
$ws_LTW.Range("I65").Value = 12227
$ws_LTW.Range("J65").Value = 24100
$ws_LTW.Range("K65").Value = 36681
$ws_LTW.Range("L65").Value = 72300
$ws_LTW.Range("M65").Value = -33561
$ws_LTW.Range("N65").Value = -78540

# Row 132 (LTW)
$ws_LTW.Range("H132").Value = 7786.9546
$ws_LTW.Range("I132").Value = 3512.2222
$ws_LTW.Range("J132").Value = 10746.385
$ws_LTW.Range("K132").Value = 10536.6666
$ws_LTW.Range("L132").Value = 32239.155
$ws_LTW.Range("M132").Value = -8006.6666
$ws_LTW.Range("N132").Value = -37299.155

# Row 140 (LTW)
$ws_LTW.Range("H140").Value = 70435.48
$ws_LTW.Range("J140").Value = 70435.48
$ws_LTW.Range("L140").Value = 70435.48
$ws_LTW.Range("N140").Value = -80795.48

$ws_WVR = $wb.Worksheets.Item("WVR")
# Row 46 (WVR)
$ws_WVR.Range("H46").Value = 49682.168
$ws_WVR.Range("J46").Value = 49682.168
$ws_WVR.Range("L46").Value = 49682.168
$ws_WVR.Range("N46").Value = -50144.168

# Row 113 (WVR)
$ws_WVR.Range("H113").Value = 5262.65
$ws_WVR.Range("J113").Value = 242.55556
$ws_WVR.Range("L113").Value = 727.66668
$ws_WVR.Range("N113").Value = -5067.66668

# Row 123 (WVR)
$ws_WVR.Range("H123").Value = 35346.152
$ws_WVR.Range("J123").Value = 35346.152
$ws_WVR.Range("L123").Value = 35346.152
$ws_WVR.Range("N123").Value = -45146.152

# Row 134 (WVR)
$ws_WVR.Range("H134").Value = 49682.168
$ws_WVR.Range("J134").Value = 49682.168
$ws_WVR.Range("L134").Value = 149046.504
$ws_WVR.Range("N134").Value = -154116.504

# Row 136 (WVR)
$ws_WVR.Range("H136").Value = 4716.213
$ws_WVR.Range("I136").Value = 4423.7334
$ws_WVR.Range("J136").Value = 5232.353
$ws_WVR.Range("K136").Value = 13271.2002
$ws_WVR.Range("L136").Value = 15697.059
$ws_WVR.Range("M136").Value = -10721.2002
$ws_WVR.Range("N136").Value = -20797.059
